# Adding objective IQ functionality
# - Insert a new "ParentIndex" column (G) before the existing Filename
#   column, shifting Filename/Response1/Response2 one column to the right.
# - Fill the new column with a constant value of 1 for each data row.
# - Remove the trailing blank rows (6-9) that only carried stray style info.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column at G, pushing Filename/Response1/Response2 -> H/I/J
$ws.Columns.Item(7).Insert()

# Header + values for the new ParentIndex column
$ws.Range("G1").Value = "ParentIndex"
$ws.Range("G2:G5").Value = 1

# Approximate the original bestFit width for the new column (closest
# value reachable through the ColumnWidth property on this host).
$ws.Columns.Item(7).ColumnWidth = 11

# Drop the now-unused blank rows 6-9
$ws.Range("A6:A9").EntireRow.Delete()
